$d = $word.ActiveDocument

# Locate the paragraph that holds the sentence we need to extend.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "This is a Microsoft word document.*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$pr = $target.Range
# Exclude the trailing paragraph mark from the range.
$pr.End = $pr.End - 1

$insertStart = $pr.End
$addition1 = " ("
$addition2 = "Changed main"
$addition3 = ")"

# Append all the new text first (this lands in a single run).
$pr.Collapse(0)  # wdCollapseEnd
$pr.InsertAfter($addition1 + $addition2 + $addition3)

$boundary1 = $insertStart
$boundary2 = $insertStart + $addition1.Length
$boundary3 = $boundary2 + $addition2.Length

# Splitting a run by planting and immediately removing a zero-length
# bookmark at a boundary forces the engine to break the text node in two
# without leaving any residual run formatting behind. Doing this from the
# right-most boundary to the left-most keeps each split crisp (matches
# how Word itself would break runs while typing left-to-right then
# normalizing).
$d.Bookmarks.Add("ironSplit3", $d.Range($boundary3, $boundary3))
$d.Bookmarks("ironSplit3").Delete()

$d.Bookmarks.Add("ironSplit2", $d.Range($boundary2, $boundary2))
$d.Bookmarks("ironSplit2").Delete()

$d.Bookmarks.Add("ironSplit1", $d.Range($boundary1, $boundary1))
$d.Bookmarks("ironSplit1").Delete()
